$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "delM%"
$ws.Range("D1").Value = "rel M"
$ws.Range("E1").Value = "rel CBF"

$ws.Range("C2").Value = 0
$ws.Range("C3").Value = -6.57
$ws.Range("C4").Value = -11.01
$ws.Range("C5").Value = -14.96
$ws.Range("C6").Value = -17.69

$ws.Range("D2").Formula = "=1+C2/100"
$ws.Range("D3:D6").Formula = "=1+C3/100"

$ws.Range("E2").Formula = "=D2*`$B`$2/B2"
$ws.Range("E3:E6").Formula = "=D3*`$B`$2/B3"

$ws.Range("H2:J5").FormulaArray = "=LINEST(B2:B6,A2:A6,TRUE,TRUE)"
